$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '''243.40'
$ws.Cells.Item(2, 7).Value = '''6'

# Row 3
$ws.Cells.Item(3, 4).Value = '''23.09'
$ws.Cells.Item(3, 7).Value = '''6'

# Row 4
$ws.Cells.Item(4, 4).Value = '''5.410'
$ws.Cells.Item(4, 7).Value = '''6'

# Row 5
$ws.Cells.Item(5, 4).Value = '''0.05990'
$ws.Cells.Item(5, 7).Value = '''6'

# Row 6
$ws.Cells.Item(6, 7).Value = '''6'

# Row 7
$ws.Cells.Item(7, 4).Value = '''6.502'
$ws.Cells.Item(7, 7).Value = '''6'

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.8107'
$ws.Cells.Item(8, 7).Value = '''6'

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.9220'
$ws.Cells.Item(9, 7).Value = '''6'

# Row 10
$ws.Cells.Item(10, 2).Value = 'WazirX'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(10, 4).Value = '''0.1425'
$ws.Cells.Item(10, 5).Value = '9WazirXWRX'
$ws.Cells.Item(10, 7).Value = '''6'

# Row 11
$ws.Cells.Item(11, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(11, 4).Value = '''0.07383'
$ws.Cells.Item(11, 5).Value = '10MandalaExchangeTokenMDX'
$ws.Cells.Item(11, 7).Value = '''6'

# Row 12
$ws.Cells.Item(12, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(12, 4).Value = '''0.03287'
$ws.Cells.Item(12, 5).Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Cells.Item(12, 7).Value = '''6'

# Row 13
$ws.Cells.Item(13, 2).Value = 'BitrueCoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(13, 4).Value = '''0.03071'
$ws.Cells.Item(13, 5).Value = '12BitrueCoinBTR'
$ws.Cells.Item(13, 7).Value = '''6'

# Row 14
$ws.Cells.Item(14, 2).Value = 'BitMartToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(14, 4).Value = '''0.09359'
$ws.Cells.Item(14, 5).Value = '13BitMartTokenBMX'
$ws.Cells.Item(14, 7).Value = '''6'

# Row 15
$ws.Cells.Item(15, 2).Value = 'MCDex'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Cells.Item(15, 4).Value = '''3.856'
$ws.Cells.Item(15, 5).Value = '14MCDexMCB'
$ws.Cells.Item(15, 7).Value = '''6'

# Row 16
$ws.Cells.Item(16, 2).Value = 'BitForexToken'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(16, 4).Value = '''0.001573'
$ws.Cells.Item(16, 5).Value = '15BitForexTokenBF'
$ws.Cells.Item(16, 7).Value = '''6'

# Row 17
$ws.Cells.Item(17, 2).Value = 'CoinExToken'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Cells.Item(17, 4).Value = '''0.04710'
$ws.Cells.Item(17, 5).Value = '16CoinExTokenCET'
$ws.Cells.Item(17, 7).Value = '''6'

# Row 18
$ws.Cells.Item(18, 2).Value = 'TigerCash'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(18, 4).Value = '''0.005870'
$ws.Cells.Item(18, 5).Value = '17TigerCashTCH'
$ws.Cells.Item(18, 7).Value = '''6'

# Row 19
$ws.Cells.Item(19, 2).Value = 'BitKan'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Cells.Item(19, 4).Value = '''0.001266'
$ws.Cells.Item(19, 5).Value = '18BitKanKAN'
$ws.Cells.Item(19, 7).Value = '''6'

# Row 20
$ws.Cells.Item(20, 2).Value = 'HotbitToken'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Cells.Item(20, 4).Value = '''0.004875'
$ws.Cells.Item(20, 5).Value = '19HotbitTokenHTB'
$ws.Cells.Item(20, 7).Value = '''6'

# Row 21
$ws.Cells.Item(21, 2).Value = 'NitroEx'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Cells.Item(21, 4).Value = '''0.00006800'
$ws.Cells.Item(21, 5).Value = '20NitroExNTX'
$ws.Cells.Item(21, 7).Value = '''6'

# Row 22
$ws.Cells.Item(22, 2).Value = 'LEO'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(22, 4).Value = '''3.575'
$ws.Cells.Item(22, 5).Value = '21LEOLEO'
$ws.Cells.Item(22, 7).Value = '''6'

# Row 23
$ws.Cells.Item(23, 2).Value = 'BTSEToken'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(23, 4).Value = '''2.133'
$ws.Cells.Item(23, 5).Value = '22BTSETokenBTSE'
$ws.Cells.Item(23, 7).Value = '''6'

# Row 24
$ws.Cells.Item(24, 2).Value = 'One'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(24, 4).Value = '''0.01115'
$ws.Cells.Item(24, 5).Value = '23OneONEBestin24h'
$ws.Cells.Item(24, 7).Value = '''6'

# Row 25
$ws.Cells.Item(25, 4).Value = '''0.3237'
$ws.Cells.Item(25, 7).Value = '''6'

# Row 26
$ws.Cells.Item(26, 4).Value = '''0.1329'
$ws.Cells.Item(26, 7).Value = '''6'

# Row 27
$ws.Cells.Item(27, 4).Value = '''0.0002340'
$ws.Cells.Item(27, 7).Value = '''6'

# Row 28
$ws.Cells.Item(28, 7).Value = '''6'

# Row 29
$ws.Cells.Item(29, 7).Value = '''6'

# Row 30
$ws.Cells.Item(30, 7).Value = '''6'

# Row 31
$ws.Cells.Item(31, 7).Value = '''6'

# Row 32
$ws.Cells.Item(32, 7).Value = '''6'

# Row 33
$ws.Cells.Item(33, 7).Value = '''6'

# Row 34
$ws.Cells.Item(34, 7).Value = '''6'

# Row 35
$ws.Cells.Item(35, 7).Value = '''6'

# Row 36
$ws.Cells.Item(36, 7).Value = '''6'

# Row 37
$ws.Cells.Item(37, 7).Value = '''6'

# Row 38
$ws.Cells.Item(38, 7).Value = '''6'

# Row 39
$ws.Cells.Item(39, 7).Value = '''6'

# Row 40
$ws.Cells.Item(40, 4).Value = '''0.03970'
$ws.Cells.Item(40, 7).Value = '''6'

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.006380'
$ws.Cells.Item(41, 7).Value = '''6'

# Row 42
$ws.Cells.Item(42, 4).Value = '''0.004300'
$ws.Cells.Item(42, 7).Value = '''6'

# Row 43
$ws.Cells.Item(43, 7).Value = '''6'

# Row 44
$ws.Cells.Item(44, 4).Value = '''0.008302'
$ws.Cells.Item(44, 7).Value = '''6'

# Row 45
$ws.Cells.Item(45, 4).Value = '''0.00005083'
$ws.Cells.Item(45, 7).Value = '''6'

# Row 46
$ws.Cells.Item(46, 7).Value = '''6'

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.6701'
$ws.Cells.Item(47, 7).Value = '''6'

# Row 48
$ws.Cells.Item(48, 4).Value = '''0.002318'
$ws.Cells.Item(48, 7).Value = '''6'

# Row 49
$ws.Cells.Item(49, 4).Value = '''0.00002100'
$ws.Cells.Item(49, 7).Value = '''6'

# Row 50
$ws.Cells.Item(50, 4).Value = '''0.0002000'
$ws.Cells.Item(50, 7).Value = '''6'

# Row 51
$ws.Cells.Item(51, 7).Value = '''6'
